# Adding sex selection to weight for age plot example.
#
# The "survey" sheet gets two new rows inserted into the custom_template
# example block (between "begin screen" and "end screen"):
#   - a new "note" row introducing the weight-for-age plot data
#   - a new "select_one sexes" row asking for the respondent's sex
# The existing "age" and "weight" prompts get their labels clarified, and
# the age prompt gets a constraint_message.
# The "choices" sheet gets a new "sexes" choice list with male/female options.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Insert a blank row above the "age" prompt (currently row 46) for the new
# introductory note, and a blank row above "end screen" (row 49, once the
# first insert has shifted everything down by one) for the new sex prompt.
$ws.Range("A46").EntireRow.Insert()
$ws.Range("A49").EntireRow.Insert()

# New row 46: introductory note for the plot.
$ws.Range("A46").Value = "note"
$ws.Range("D46").Value = "The following data will be used to generate a weight for age plot."

# Existing "age" prompt, now shifted down to row 47: clarify label & add a
# constraint message.
$ws.Range("D47").Value = "Enter age (in years):"
$ws.Range("E47").Value = "Must be less than 20."

# Existing "weight" prompt, now shifted down to row 48: clarify label.
$ws.Range("D48").Value = "Enter weight (in lbs):"

# New row 49: ask for sex using the new "sexes" choice list.
$ws.Range("A49").Value = "select_one sexes"
$ws.Range("C49").Value = "sex"
$ws.Range("D49").Value = "Enter sex:"

# Add the new "sexes" choice list (male/female) to the "choices" sheet.
$wsChoices = $wb.Worksheets.Item("choices")
$wsChoices.Range("A23").Value = "sexes"
$wsChoices.Range("B23").Value = "male"
$wsChoices.Range("C23").Value = "male"
$wsChoices.Range("A24").Value = "sexes"
$wsChoices.Range("B24").Value = "female"
$wsChoices.Range("C24").Value = "female"
